$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill C24:F24 with the value 5 (row 21 - "the event ... happened")
$ws.Range("C24:F24").Value = 5

# Move the active selection in the bottom-right frozen pane to G24
$ws.Range("G24").Select()
